$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.397.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.17%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.793.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.77%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.28%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5348"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.93%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3763"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.70%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07489"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.07%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.06%  "

# Row 11
$ws.Range("E11").Value = "  -1.87%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.18%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.79%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.128"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.00%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.806.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.25%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.282"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.89%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.89%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001060"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.91%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06490"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.43%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.22%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.71%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.963"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.442.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.20%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.082"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.08%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.91%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.64%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.012.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.11%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.299"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.91%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.46%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.098"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.82%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1047"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.39%  "

# Row 33
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.660"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.586"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.29%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06552"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.52%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2256"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.62%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02292"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.36%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.994"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.38%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.516"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.14%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.451"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.22%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6174"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.98%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.192"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.03%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.15%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.51%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.681"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.17%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5784"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.79%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.70%  "

# Row 49
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.191"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.68%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.938"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.42%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06863"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.67%  "
